$d = $word.ActiveDocument

$d.Content.Find.Execute("94-4=90", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2) | Out-Null
$d.Content.Find.Execute("35+33=68", $true, $false, $false, $false, $false, $true, 1, $false, "42+10=52", 2) | Out-Null
$d.Content.Find.Execute("56-50=6", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=83", 2) | Out-Null
$d.Content.Find.Execute("78+19=97", $true, $false, $false, $false, $false, $true, 1, $false, "17+81=98", 2) | Out-Null
$d.Content.Find.Execute("53-25=28", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=6", 2) | Out-Null
$d.Content.Find.Execute("6+10=16", $true, $false, $false, $false, $false, $true, 1, $false, "71-45=26", 2) | Out-Null
$d.Content.Find.Execute("55+20=75", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2) | Out-Null
$d.Content.Find.Execute("8+39=47", $true, $false, $false, $false, $false, $true, 1, $false, "73+18=91", 2) | Out-Null
$d.Content.Find.Execute("14-10=4", $true, $false, $false, $false, $false, $true, 1, $false, "10+50=60", 2) | Out-Null
$d.Content.Find.Execute("17+62=79", $true, $false, $false, $false, $false, $true, 1, $false, "89-16=73", 2) | Out-Null
$d.Content.Find.Execute("77+14=91", $true, $false, $false, $false, $false, $true, 1, $false, "51+31=82", 2) | Out-Null
$d.Content.Find.Execute("11+25=36", $true, $false, $false, $false, $false, $true, 1, $false, "97-10=87", 2) | Out-Null
$d.Content.Find.Execute("23+2=25", $true, $false, $false, $false, $false, $true, 1, $false, "79-7=72", 2) | Out-Null
$d.Content.Find.Execute("41+12=53", $true, $false, $false, $false, $false, $true, 1, $false, "39+18=57", 2) | Out-Null
$d.Content.Find.Execute("90-28=62", $true, $false, $false, $false, $false, $true, 1, $false, "28+59=87", 2) | Out-Null
$d.Content.Find.Execute("99-15=84", $true, $false, $false, $false, $false, $true, 1, $false, "11+34=45", 2) | Out-Null
$d.Content.Find.Execute("37+8=45", $true, $false, $false, $false, $false, $true, 1, $false, "89-12=77", 2) | Out-Null
$d.Content.Find.Execute("24+21=45", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=79", 2) | Out-Null
$d.Content.Find.Execute("45-44=1", $true, $false, $false, $false, $false, $true, 1, $false, "80-77=3", 2) | Out-Null
$d.Content.Find.Execute("36+6=42", $true, $false, $false, $false, $false, $true, 1, $false, "95-74=21", 2) | Out-Null
$d.Content.Find.Execute("30+30=60", $true, $false, $false, $false, $false, $true, 1, $false, "61-35=26", 2) | Out-Null
$d.Content.Find.Execute("20+75=95", $true, $false, $false, $false, $false, $true, 1, $false, "73+6=79", 2) | Out-Null
$d.Content.Find.Execute("8+18=26", $true, $false, $false, $false, $false, $true, 1, $false, "51+7=58", 2) | Out-Null
$d.Content.Find.Execute("41+8=49", $true, $false, $false, $false, $false, $true, 1, $false, "62-60=2", 2) | Out-Null
$d.Content.Find.Execute("65+34=99", $true, $false, $false, $false, $false, $true, 1, $false, "26+52=78", 2) | Out-Null
$d.Content.Find.Execute("16+23=39", $true, $false, $false, $false, $false, $true, 1, $false, "51-22=29", 2) | Out-Null
$d.Content.Find.Execute("60+1=61", $true, $false, $false, $false, $false, $true, 1, $false, "48-1=47", 2) | Out-Null
$d.Content.Find.Execute("17+66=83", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("46-32=14", $true, $false, $false, $false, $false, $true, 1, $false, "33+5=38", 2) | Out-Null
$d.Content.Find.Execute("76-26=50", $true, $false, $false, $false, $false, $true, 1, $false, "88-55=33", 2) | Out-Null
$d.Content.Find.Execute("38+58=96", $true, $false, $false, $false, $false, $true, 1, $false, "47-33=14", 2) | Out-Null
$d.Content.Find.Execute("44-14=30", $true, $false, $false, $false, $false, $true, 1, $false, "28+34=62", 2) | Out-Null
$d.Content.Find.Execute("3+6=9", $true, $false, $false, $false, $false, $true, 1, $false, "25+74=99", 2) | Out-Null
$d.Content.Find.Execute("90-87=3", $true, $false, $false, $false, $false, $true, 1, $false, "42-30=12", 2) | Out-Null
$d.Content.Find.Execute("88-40=48", $true, $false, $false, $false, $false, $true, 1, $false, "80-5=75", 2) | Out-Null
$d.Content.Find.Execute("67+21=88", $true, $false, $false, $false, $false, $true, 1, $false, "70-56=14", 2) | Out-Null
$d.Content.Find.Execute("31+18=49", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=38", 2) | Out-Null
$d.Content.Find.Execute("74-57=17", $true, $false, $false, $false, $false, $true, 1, $false, "14-8=6", 2) | Out-Null
$d.Content.Find.Execute("18+45=63", $true, $false, $false, $false, $false, $true, 1, $false, "56-22=34", 2) | Out-Null
$d.Content.Find.Execute("72-57=15", $true, $false, $false, $false, $false, $true, 1, $false, "98-94=4", 2) | Out-Null
$d.Content.Find.Execute("82-57=25", $true, $false, $false, $false, $false, $true, 1, $false, "69-35=34", 2) | Out-Null
$d.Content.Find.Execute("41-39=2", $true, $false, $false, $false, $false, $true, 1, $false, "54-34=20", 2) | Out-Null
$d.Content.Find.Execute("69-24=45", $true, $false, $false, $false, $false, $true, 1, $false, "68-36=32", 2) | Out-Null
$d.Content.Find.Execute("91+5=96", $true, $false, $false, $false, $false, $true, 1, $false, "47-26=21", 2) | Out-Null
$d.Content.Find.Execute("31-25=6", $true, $false, $false, $false, $false, $true, 1, $false, "65+30=95", 2) | Out-Null
$d.Content.Find.Execute("43+41=84", $true, $false, $false, $false, $false, $true, 1, $false, "66-58=8", 2) | Out-Null
$d.Content.Find.Execute("94-50=44", $true, $false, $false, $false, $false, $true, 1, $false, "96-53=43", 2) | Out-Null
$d.Content.Find.Execute("53+39=92", $true, $false, $false, $false, $false, $true, 1, $false, "69-57=12", 2) | Out-Null
$d.Content.Find.Execute("44+31=75", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("91-49=42", $true, $false, $false, $false, $false, $true, 1, $false, "51-46=5", 2) | Out-Null
$d.Content.Find.Execute("21-10=11", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=48", 2) | Out-Null
$d.Content.Find.Execute("61+32=93", $true, $false, $false, $false, $false, $true, 1, $false, "61+3=64", 2) | Out-Null
$d.Content.Find.Execute("72-45=27", $true, $false, $false, $false, $false, $true, 1, $false, "84+0=84", 2) | Out-Null
$d.Content.Find.Execute("60+0=60", $true, $false, $false, $false, $false, $true, 1, $false, "36+15=51", 2) | Out-Null
$d.Content.Find.Execute("92-49=43", $true, $false, $false, $false, $false, $true, 1, $false, "94-79=15", 2) | Out-Null
$d.Content.Find.Execute("27+8=35", $true, $false, $false, $false, $false, $true, 1, $false, "51+29=80", 2) | Out-Null
$d.Content.Find.Execute("72-12=60", $true, $false, $false, $false, $false, $true, 1, $false, "65-0=65", 2) | Out-Null
$d.Content.Find.Execute("91-20=71", $true, $false, $false, $false, $false, $true, 1, $false, "39+7=46", 2) | Out-Null
$d.Content.Find.Execute("0+66=66", $true, $false, $false, $false, $false, $true, 1, $false, "47-9=38", 2) | Out-Null
$d.Content.Find.Execute("78-50=28", $true, $false, $false, $false, $false, $true, 1, $false, "55+16=71", 2) | Out-Null
$d.Content.Find.Execute("44+45=89", $true, $false, $false, $false, $false, $true, 1, $false, "58+21=79", 2) | Out-Null
$d.Content.Find.Execute("26+56=82", $true, $false, $false, $false, $false, $true, 1, $false, "60+37=97", 2) | Out-Null
$d.Content.Find.Execute("63-43=20", $true, $false, $false, $false, $false, $true, 1, $false, "38-8=30", 2) | Out-Null
$d.Content.Find.Execute("76+16=92", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 2) | Out-Null
$d.Content.Find.Execute("9+62=71", $true, $false, $false, $false, $false, $true, 1, $false, "90-85=5", 2) | Out-Null
$d.Content.Find.Execute("12+38=50", $true, $false, $false, $false, $false, $true, 1, $false, "72+21=93", 2) | Out-Null
$d.Content.Find.Execute("80-46=34", $true, $false, $false, $false, $false, $true, 1, $false, "0+43=43", 2) | Out-Null
$d.Content.Find.Execute("96-62=34", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=41", 2) | Out-Null
$d.Content.Find.Execute("65-6=59", $true, $false, $false, $false, $false, $true, 1, $false, "79-5=74", 2) | Out-Null
$d.Content.Find.Execute("11-4=7", $true, $false, $false, $false, $false, $true, 1, $false, "97-11=86", 2) | Out-Null
$d.Content.Find.Execute("38-25=13", $true, $false, $false, $false, $false, $true, 1, $false, "80+15=95", 2) | Out-Null
$d.Content.Find.Execute("35+14=49", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=76", 2) | Out-Null
$d.Content.Find.Execute("35-30=5", $true, $false, $false, $false, $false, $true, 1, $false, "90-2=88", 2) | Out-Null
$d.Content.Find.Execute("37-8=29", $true, $false, $false, $false, $false, $true, 1, $false, "89-14=75", 2) | Out-Null
$d.Content.Find.Execute("4+75=79", $true, $false, $false, $false, $false, $true, 1, $false, "86-2=84", 2) | Out-Null
$d.Content.Find.Execute("71+26=97", $true, $false, $false, $false, $false, $true, 1, $false, "22+25=47", 2) | Out-Null
$d.Content.Find.Execute("75-11=64", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=93", 2) | Out-Null
$d.Content.Find.Execute("70-21=49", $true, $false, $false, $false, $false, $true, 1, $false, "49+4=53", 2) | Out-Null
$d.Content.Find.Execute("1+47=48", $true, $false, $false, $false, $false, $true, 1, $false, "15+14=29", 2) | Out-Null
$d.Content.Find.Execute("16+12=28", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=6", 2) | Out-Null
$d.Content.Find.Execute("46+32=78", $true, $false, $false, $false, $false, $true, 1, $false, "23+16=39", 2) | Out-Null
$d.Content.Find.Execute("54+16=70", $true, $false, $false, $false, $false, $true, 1, $false, "77-50=27", 2) | Out-Null
$d.Content.Find.Execute("9+37=46", $true, $false, $false, $false, $false, $true, 1, $false, "96-85=11", 2) | Out-Null
$d.Content.Find.Execute("32+16=48", $true, $false, $false, $false, $false, $true, 1, $false, "94-38=56", 2) | Out-Null
$d.Content.Find.Execute("22+76=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+29=31", 2) | Out-Null
$d.Content.Find.Execute("74-23=51", $true, $false, $false, $false, $false, $true, 1, $false, "37-12=25", 2) | Out-Null
$d.Content.Find.Execute("97-91=6", $true, $false, $false, $false, $false, $true, 1, $false, "83+9=92", 2) | Out-Null
$d.Content.Find.Execute("7+10=17", $true, $false, $false, $false, $false, $true, 1, $false, "72+8=80", 2) | Out-Null
$d.Content.Find.Execute("83-42=41", $true, $false, $false, $false, $false, $true, 1, $false, "0+67=67", 2) | Out-Null
$d.Content.Find.Execute("81+18=99", $true, $false, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null
$d.Content.Find.Execute("12-7=5", $true, $false, $false, $false, $false, $true, 1, $false, "90-89=1", 2) | Out-Null
$d.Content.Find.Execute("61-59=2", $true, $false, $false, $false, $false, $true, 1, $false, "53-38=15", 2) | Out-Null
$d.Content.Find.Execute("71-54=17", $true, $false, $false, $false, $false, $true, 1, $false, "87-54=33", 2) | Out-Null
$d.Content.Find.Execute("54-44=10", $true, $false, $false, $false, $false, $true, 1, $false, "56+23=79", 2) | Out-Null
$d.Content.Find.Execute("44+16=60", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=82", 2) | Out-Null
$d.Content.Find.Execute("73+2=75", $true, $false, $false, $false, $false, $true, 1, $false, "27-26=1", 2) | Out-Null
$d.Content.Find.Execute("13+85=98", $true, $false, $false, $false, $false, $true, 1, $false, "60-42=18", 2) | Out-Null
$d.Content.Find.Execute("25+55=80", $true, $false, $false, $false, $false, $true, 1, $false, "48+46=94", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $false, $false, $false, $false, $true, 1, $false, "7+52=59", 2) | Out-Null
$d.Content.Find.Execute("75-41=34", $true, $false, $false, $false, $false, $true, 1, $false, "1+73=74", 2) | Out-Null

Write-Output "Done: replaced 100 values"
